$d = $word.ActiveDocument

# Locate the paragraph holding objective "4 ." — the new objectives 5 and 6
# must be inserted right after it, before the trailing blank paragraphs.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Avoir un sens de travail*envergure*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq 0) {
    throw "Could not find paragraph 4 (objectif évolutif) to anchor the new items after."
}

# Insert objective 5 as a brand-new paragraph right after objective 4.
$anchor = $d.Paragraphs($targetIndex).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$p5 = $d.Paragraphs($targetIndex + 1)
$p5.Range.InsertBefore("5. Etre parmi les meilleurs de la promotion")

# Insert objective 6 as a brand-new paragraph right after objective 5.
$anchor2 = $d.Paragraphs($targetIndex + 1).Range
$anchor2.Collapse(0)
$anchor2.InsertParagraphAfter()
$p6 = $d.Paragraphs($targetIndex + 2)
$p6.Range.InsertBefore("6. Trouver de l’emploi après la certification")
